$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Customer_ID values (column B) for rows 2-11 on Sheet1.
# Force text formatting first so the numeric-looking IDs are stored as
# shared strings (matching the workbook's original text representation),
# then reset the cell style so no stray NumberFormat is left behind.
$rng = $ws.Range("B2:B11")
$rng.NumberFormat = "@"

$ws.Range("B2").Value = "17704808"
$ws.Range("B3").Value = "17704809"
$ws.Range("B4").Value = "17704810"
$ws.Range("B5").Value = "17704811"
$ws.Range("B6").Value = "17704812"
$ws.Range("B7").Value = "17704813"
$ws.Range("B8").Value = "17704815"
$ws.Range("B9").Value = "17704816"
$ws.Range("B10").Value = "17704817"
$ws.Range("B11").Value = "17704818"

$rng.Style = "Normal"

# Update the selection range on Sheet1 from A2:C2 to A2:C6
$ws.Range("A2:C6").Select()
